# Applies the "Estadisticos Segundo Parcial 23 Mayo" update:
# - "2o Parcial": fills in real Aprobados/Reprobados/Promedio/Blancos stats
#   (previously every student was counted as "Blancos" / no-show).
# - "Final": recomputed combined Aprobados/Reprobados/Promedio now that the
#   Segundo Parcial results are known.
$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("2o Parcial")
$ws2.Range("E2").Value = 22
$ws2.Range("F2").Value = 9
$ws2.Range("G2").Value = 71
$ws2.Range("H2").Value = 29
$ws2.Range("I2").Value = 6.9
$ws2.Range("J2").Value = 0
$ws2.Range("K2").Value = 0

$ws2.Range("E3").Value = 38
$ws2.Range("F3").Value = 11
$ws2.Range("G3").Value = 77.59999999999999
$ws2.Range("H3").Value = 22.4
$ws2.Range("I3").Value = 7.4
$ws2.Range("J3").Value = 0
$ws2.Range("K3").Value = 0

$ws2.Range("E4").Value = 60
$ws2.Range("F4").Value = 20
$ws2.Range("G4").Value = 75
$ws2.Range("H4").Value = 25
$ws2.Range("I4").Value = 7.2
$ws2.Range("J4").Value = 0
$ws2.Range("K4").Value = 0

$ws2.Range("E7").Value = 34
$ws2.Range("F7").Value = 2
$ws2.Range("G7").Value = 94.40000000000001
$ws2.Range("H7").Value = 5.6
$ws2.Range("I7").Value = 8.199999999999999
$ws2.Range("J7").Value = 0
$ws2.Range("K7").Value = 0

$ws2.Range("E8").Value = 39
$ws2.Range("F8").Value = 0
$ws2.Range("G8").Value = 100
$ws2.Range("H8").Value = 0
$ws2.Range("I8").Value = 9.199999999999999
$ws2.Range("J8").Value = 0
$ws2.Range("K8").Value = 0

$ws2.Range("E9").Value = 39
$ws2.Range("F9").Value = 0
$ws2.Range("G9").Value = 100
$ws2.Range("H9").Value = 0
$ws2.Range("I9").Value = 8.6
$ws2.Range("J9").Value = 0
$ws2.Range("K9").Value = 0

$ws2.Range("E10").Value = 37
$ws2.Range("F10").Value = 1
$ws2.Range("G10").Value = 97.40000000000001
$ws2.Range("H10").Value = 2.6
$ws2.Range("I10").Value = 9.300000000000001
$ws2.Range("J10").Value = 0
$ws2.Range("K10").Value = 0

$ws2.Range("E11").Value = 27
$ws2.Range("F11").Value = 1
$ws2.Range("G11").Value = 96.40000000000001
$ws2.Range("H11").Value = 3.6
$ws2.Range("I11").Value = 7.9
$ws2.Range("J11").Value = 0
$ws2.Range("K11").Value = 0

$ws2.Range("E12").Value = 23
$ws2.Range("F12").Value = 0
$ws2.Range("G12").Value = 100
$ws2.Range("H12").Value = 0
$ws2.Range("I12").Value = 8
$ws2.Range("J12").Value = 0
$ws2.Range("K12").Value = 0

$ws2.Range("E13").Value = 28
$ws2.Range("F13").Value = 2
$ws2.Range("G13").Value = 93.3
$ws2.Range("H13").Value = 6.7
$ws2.Range("I13").Value = 8.5
$ws2.Range("J13").Value = 0
$ws2.Range("K13").Value = 0

$ws2.Range("E14").Value = 9
$ws2.Range("F14").Value = 2
$ws2.Range("G14").Value = 81.8
$ws2.Range("H14").Value = 18.2
$ws2.Range("I14").Value = 8.4
$ws2.Range("J14").Value = 0
$ws2.Range("K14").Value = 0

$ws2.Range("E15").Value = 236
$ws2.Range("F15").Value = 8
$ws2.Range("G15").Value = 96.7
$ws2.Range("H15").Value = 3.3
$ws2.Range("I15").Value = 8.5
$ws2.Range("J15").Value = 0
$ws2.Range("K15").Value = 0

$ws2.Range("E16").Value = 296
$ws2.Range("F16").Value = 63
$ws2.Range("G16").Value = 82.5
$ws2.Range("H16").Value = 17.5
$ws2.Range("I16").Value = 7.5
$ws2.Range("J16").Value = 35
$ws2.Range("K16").Value = 9.699999999999999

$ws3 = $wb.Worksheets.Item("Final")
$ws3.Range("E2").Value = 22
$ws3.Range("F2").Value = 9
$ws3.Range("G2").Value = 71
$ws3.Range("H2").Value = 29
$ws3.Range("I2").Value = 6.7

$ws3.Range("I3").Value = 7

$ws3.Range("E4").Value = 60
$ws3.Range("F4").Value = 20
$ws3.Range("G4").Value = 75
$ws3.Range("H4").Value = 25
$ws3.Range("I4").Value = 6.8

$ws3.Range("I7").Value = 8.1

$ws3.Range("I8").Value = 9.1

$ws3.Range("I9").Value = 8.6

$ws3.Range("I10").Value = 9.199999999999999

$ws3.Range("E11").Value = 27
$ws3.Range("F11").Value = 1
$ws3.Range("G11").Value = 96.40000000000001
$ws3.Range("H11").Value = 3.6
$ws3.Range("I11").Value = 7.6

$ws3.Range("E12").Value = 23
$ws3.Range("F12").Value = 0
$ws3.Range("G12").Value = 100
$ws3.Range("H12").Value = 0
$ws3.Range("I12").Value = 7.8

$ws3.Range("E13").Value = 28
$ws3.Range("F13").Value = 2
$ws3.Range("G13").Value = 93.3
$ws3.Range("H13").Value = 6.7
$ws3.Range("I13").Value = 8.300000000000001

$ws3.Range("E14").Value = 9
$ws3.Range("F14").Value = 2
$ws3.Range("G14").Value = 81.8
$ws3.Range("H14").Value = 18.2
$ws3.Range("I14").Value = 8.300000000000001

$ws3.Range("E15").Value = 236
$ws3.Range("F15").Value = 8
$ws3.Range("G15").Value = 96.7
$ws3.Range("H15").Value = 3.3
$ws3.Range("I15").Value = 8.4

$ws3.Range("E16").Value = 330
$ws3.Range("F16").Value = 29
$ws3.Range("G16").Value = 91.90000000000001
$ws3.Range("H16").Value = 8.1
$ws3.Range("I16").Value = 8
